$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add two new trailing data points (X3, Y3) ---
$ws.Cells.Item(3, 24).Value = -1.3299870000000169
$ws.Cells.Item(3, 25).Value = "Down"

# --- Row 4: brand new data row ---
$ws.Cells.Item(4, 1).Value = 42633.888368055559
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = "Neutral"
$ws.Cells.Item(4, 4).Value = 10
$ws.Cells.Item(4, 5).Value = 25422
$ws.Cells.Item(4, 6).Value = 2956
$ws.Cells.Item(4, 7).Value = 56
$ws.Cells.Item(4, 8).Value = 42
$ws.Cells.Item(4, 9).Value = 82
$ws.Cells.Item(4, 10).Value = 17
$ws.Cells.Item(4, 11).Value = 20655
$ws.Cells.Item(4, 12).Value = 352
$ws.Cells.Item(4, 13).Value = 267
$ws.Cells.Item(4, 14).Value = 42
$ws.Cells.Item(4, 15).Value = 9
$ws.Cells.Item(4, 16).Value = "Bag"
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 1.76
$ws.Cells.Item(4, 19).Value = 0.1055
$ws.Cells.Item(4, 19).NumberFormat = "0.00%"
$ws.Cells.Item(4, 20).Value = -6.67
$ws.Cells.Item(4, 21).Value = 5.83
$ws.Cells.Item(4, 22).Value = "N/A"
$ws.Cells.Item(4, 23).Value = 0
